# Applies the commit "added new blog entry and accompanying images":
# 5 text edits, each splitting / rewording an existing run into multiple runs
# while preserving the original run formatting (rPr).

function Split-Run {
    param($FindText, $Parts)

    $d = $word.ActiveDocument
    $rng = $d.Content
    $found = $rng.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $FindText"
        return
    }
    $start = $rng.Start

    $newFull = [string]::Join("", $Parts)
    # Replace whole matched range with the new full text (collapses to one run
    # using the original run's formatting).
    $rng.Text = $newFull

    # Now force run boundaries between the requested parts by toggling a
    # formatting property (Bold on, then back off) on each part after the
    # first one. Toggling like this causes Word to split the run without
    # altering the resulting formatting.
    $pos = $start
    for ($i = 0; $i -lt $Parts.Length; $i++) {
        $len = $Parts[$i].Length
        if ($i -gt 0) {
            $subRng = $d.Range($pos, $pos + $len)
            $subRng.Bold = 1
            $subRng.Bold = 0
        }
        $pos = $pos + $len
    }
}

# 1) "...which I have written about in my blog, linked to above" ->
#    "...which I have written about in my blog, as seen above" (split into 3 runs)
Split-Run "which I have written about in my blog, linked to above" @("which I have written about in my blog, ", "as seen", " above")

# 2) " managing an artist's website on Cargo Collective. I am currently working on " ->
#    " managing a landscape designer's website on Cargo Collective. I am currently working on " (split into 3 runs)
Split-Run " managing an artist’s website on Cargo Collective. I am currently working on " @(" managing a", " landscape designer", "’s website on Cargo Collective. I am currently working on ")

# 3) "some " -> "beginner " (split into 2 runs)
Split-Run "some " @("beginner", " ")

# 4) "Java, ISTQB Certified Software Tester, " -> "Java; ISTQB Certified Software Tester, " (split into 2 runs)
Split-Run "Java, ISTQB Certified Software Tester, " @("Java;", " ISTQB Certified Software Tester, ")

# 5) "Currently, I am helping to manage an artist's website" ->
#    "Currently, I am helping to manage an landscape designer's website" (split into 3 runs)
Split-Run "Currently, I am helping to manage an artist’s website" @("Currently, I ", "am helping to manage an landscape designer", "’s website")
